# Scheduled data refresh: update cached market-price / profit figures
# across all 8 crafting-class sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 248.91667
$ws.Range("I6").Value = 98.5
$ws.Range("K6").Value = 295.5
$ws.Range("M6").Value = -183.5
$ws.Range("H18").Value = 6742
$ws.Range("I18").Value = 3850.5
$ws.Range("K18").Value = 3850.5
$ws.Range("M18").Value = -3566.5
$ws.Range("H52").Value = 750
$ws.Range("I52").Value = 500
$ws.Range("K52").Value = 1500
$ws.Range("M52").Value = -1340
$ws.Range("H58").Value = 3200
$ws.Range("H113").Value = 4005
$ws.Range("I113").Value = 4005
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 4005
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -751
$ws.Range("N113").ClearContents()
$ws.Range("H137").Value = 3782.7144
$ws.Range("I137").Value = 9750.25
$ws.Range("J137").Value = 2378.5881
$ws.Range("K137").Value = 29250.75
$ws.Range("L137").Value = 7135.7643
$ws.Range("M137").Value = -26700.75
$ws.Range("N137").Value = -12235.7643
$ws.Range("H138").Value = 3275.5
$ws.Range("I138").Value = 2584.652
$ws.Range("J138").Value = 4032.1428
$ws.Range("K138").Value = 7753.956
$ws.Range("L138").Value = 12096.4284
$ws.Range("M138").Value = -2613.956
$ws.Range("N138").Value = -22376.4284

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 282
$ws.Range("I5").Value = 282
$ws.Range("K5").Value = 282
$ws.Range("M5").Value = -170
$ws.Range("H32").Value = 239354.45
$ws.Range("I32").Value = 794.2973
$ws.Range("K32").Value = 794.2973
$ws.Range("M32").Value = -507.2973
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()
$ws.Range("H61").Value = 4830.231
$ws.Range("I61").Value = 4779.6
$ws.Range("K61").Value = 4779.6
$ws.Range("M61").Value = -4567.6
$ws.Range("H74").Value = 2167.7144
$ws.Range("I74").Value = 1786.1818
$ws.Range("K74").Value = 1786.1818
$ws.Range("M74").Value = -912.1818000000001
$ws.Range("H77").Value = 2167.7144
$ws.Range("I77").Value = 1786.1818
$ws.Range("K77").Value = 8930.909
$ws.Range("M77").Value = -4562.909
$ws.Range("H132").Value = 3250
$ws.Range("I132").Value = 3250
$ws.Range("K132").Value = 9750
$ws.Range("M132").Value = -7220
$ws.Range("H136").Value = 4830.231
$ws.Range("I136").Value = 4779.6
$ws.Range("K136").Value = 14338.8
$ws.Range("M136").Value = -11788.8

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 282
$ws.Range("I4").Value = 282
$ws.Range("K4").Value = 282
$ws.Range("M4").Value = -167
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws.Range("H134").Value = 3390
$ws.Range("I134").Value = 3173.5
$ws.Range("K134").Value = 9520.5
$ws.Range("M134").Value = -6985.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 769944.0600000001
$ws.Range("J6").Value = 10000000
$ws.Range("L6").Value = 10000000
$ws.Range("N6").Value = -10000226
$ws.Range("H7").Value = 394.17856
$ws.Range("I7").Value = 417.54544
$ws.Range("J7").Value = 379.05884
$ws.Range("K7").Value = 417.54544
$ws.Range("L7").Value = 379.05884
$ws.Range("M7").Value = -304.54544
$ws.Range("N7").Value = -605.0588399999999
$ws.Range("H41").Value = 11221.444
$ws.Range("I41").Value = 6499.75
$ws.Range("J41").Value = 14998.8
$ws.Range("K41").Value = 6499.75
$ws.Range("L41").Value = 14998.8
$ws.Range("M41").Value = -6071.75
$ws.Range("N41").Value = -15854.8
$ws.Range("H42").Value = 3999.6667
$ws.Range("I42").Value = 3999.6667
$ws.Range("K42").Value = 3999.6667
$ws.Range("M42").Value = -3406.6667
$ws.Range("H56").Value = 34522.5
$ws.Range("I56").Value = 93
$ws.Range("J56").Value = 45999
$ws.Range("K56").Value = 93
$ws.Range("L56").Value = 45999
$ws.Range("M56").Value = 752
$ws.Range("N56").Value = -47689
$ws.Range("H62").Value = 86621
$ws.Range("I62").Value = 7965.3335
$ws.Range("J62").Value = 204604.5
$ws.Range("K62").Value = 7965.3335
$ws.Range("L62").Value = 204604.5
$ws.Range("M62").Value = -7341.3335
$ws.Range("N62").Value = -205852.5
$ws.Range("H65").Value = 86621
$ws.Range("I65").Value = 7965.3335
$ws.Range("J65").Value = 204604.5
$ws.Range("K65").Value = 39826.6675
$ws.Range("L65").Value = 1023022.5
$ws.Range("M65").Value = -36706.6675
$ws.Range("N65").Value = -1029262.5
$ws.Range("H86").Value = 8555
$ws.Range("I86").Value = 8459.4
$ws.Range("K86").Value = 8459.4
$ws.Range("M86").Value = -7336.4
$ws.Range("H89").Value = 8555
$ws.Range("I89").Value = 8459.4
$ws.Range("K89").Value = 42297
$ws.Range("M89").Value = -36681

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 519
$ws.Range("I5").Value = 534.5
$ws.Range("J5").Value = 498.33334
$ws.Range("K5").Value = 1603.5
$ws.Range("L5").Value = 1495.00002
$ws.Range("M5").Value = -1491.5
$ws.Range("N5").Value = -1719.00002
$ws.Range("H11").Value = 1380.1578
$ws.Range("I11").Value = 513.06665
$ws.Range("K11").Value = 1539.19995
$ws.Range("M11").Value = -1399.19995
$ws.Range("H22").Value = 1629
$ws.Range("J22").Value = 1536.25
$ws.Range("L22").Value = 4608.75
$ws.Range("N22").Value = -4946.75
$ws.Range("H27").Value = 1629
$ws.Range("J27").Value = 1536.25
$ws.Range("L27").Value = 4608.75
$ws.Range("N27").Value = -4812.75
$ws.Range("H113").Value = 863.63635
$ws.Range("J113").Value = 1031.3334
$ws.Range("L113").Value = 3094.0002
$ws.Range("N113").Value = -7434.0002
$ws.Range("H121").Value = 328.2
$ws.Range("I121").Value = 160.25
$ws.Range("J121").Value = 1000
$ws.Range("K121").Value = 480.75
$ws.Range("L121").Value = 3000
$ws.Range("M121").Value = 829.25
$ws.Range("N121").Value = -5620
$ws.Range("H122").Value = 1359.2
$ws.Range("J122").Value = 1542.1666
$ws.Range("L122").Value = 13879.4994
$ws.Range("N122").Value = -18779.4994
$ws.Range("H132").Value = 4455.3125
$ws.Range("I132").Value = 2873.75
$ws.Range("K132").Value = 25863.75
$ws.Range("M132").Value = -23333.75
$ws.Range("H135").Value = 519
$ws.Range("I135").Value = 534.5
$ws.Range("J135").Value = 498.33334
$ws.Range("K135").Value = 4810.5
$ws.Range("L135").Value = 4485.00006
$ws.Range("M135").Value = -2275.5
$ws.Range("N135").Value = -9555.00006

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H49").Value = 30000
$ws.Range("J49").Value = 30000
$ws.Range("L49").Value = 30000
$ws.Range("N49").Value = -30368
$ws.Range("H63").Value = 53797.8
$ws.Range("J63").Value = 54249.75
$ws.Range("L63").Value = 54249.75
$ws.Range("N63").Value = -55621.75
$ws.Range("H66").Value = 53797.8
$ws.Range("J66").Value = 54249.75
$ws.Range("L66").Value = 162749.25
$ws.Range("N66").Value = -169613.25
$ws.Range("H80").Value = 1665
$ws.Range("H83").Value = 1665
$ws.Range("H132").Value = 2076.05
$ws.Range("I132").Value = 1732.625
$ws.Range("K132").Value = 5197.875
$ws.Range("M132").Value = -2667.875

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H104").Value = 23633
$ws.Range("J104").Value = 23633
$ws.Range("L104").Value = 23633
$ws.Range("N104").Value = -30621

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H11").Value = 15995
$ws.Range("J11").Value = 15995
$ws.Range("L11").Value = 15995
$ws.Range("N11").Value = -16279
$ws.Range("H39").Value = 6300
$ws.Range("I39").Value = 3350
$ws.Range("K39").Value = 3350
$ws.Range("M39").Value = -2937
$ws.Range("H95").Value = 14336
$ws.Range("J95").Value = 14336
$ws.Range("L95").Value = 14336
$ws.Range("N95").Value = -19828
$ws.Range("H132").Value = 6646.857
$ws.Range("J132").Value = 9864.857
$ws.Range("L132").Value = 29594.571
$ws.Range("N132").Value = -34654.571
